# Register with Invalid email
# Adds a new "Invalid_Email" test-data column (F) to Sheet1 and a new data
# row (row 4) containing an invalid email test case, including its
# corresponding mailto hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column F header + existing rows' new column values -------------
$ws.Range("F1").Value = "Invalid_Email"
$ws.Range("F2").Value = "userexample.com"
$ws.Range("F3").Value = "user.example.com"

# --- New row 4 of test data ----------------------------------------------
# NOTE: column order below matters only to control the order in which new
# strings are appended to the shared-string table; it mirrors the order
# used when this test row was originally authored.
$ws.Range("A4").Value = "Lopez"
$ws.Range("B4").Value = "Maria"
$ws.Range("D4").Value = "Ravi@2024"
$ws.Range("E4").Value = "Ravi@2024"
$ws.Range("F4").Value = "@domain.com"
$ws.Range("C4").Value = "maria.lopez2@stmail.com"

# Hyperlink the new invalid e-mail address cell, same as the other Email
# column cells (C2, C3) already have.
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:maria.lopez2@stmail.com")
# Hyperlinks.Add re-applies hyperlink formatting with a brand new style;
# put it back on the same "Hyperlink" cell style the rest of the column uses.
$ws.Range("C4").Style = "Hyperlink"

# Move the active selection, matching where editing left off.
[void]$ws.Range("F11").Select()
